$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for A2:D11 (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$data = @(
    @(9, 2, 5, 5),
    @(3, 4, 5, 5),
    @(5, 4, 10, 10),
    @(1, 6, 5, 6),
    @(2, 6, 11, 12),
    @(4, 6, 17, 18),
    @(6, 6, 23, 25),
    @(7, 6, 30, 30),
    @(10, 6, 35, 36),
    @(8, 5, 5, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$wb.Save()
